# Fixed allocation of archetypes creating duplicates
# Slide 8 ("Engine (ECS)") - TextBox 4 (shape 3): remove the stray
# "1111...1" placeholder line, add a note about profiling before the
# "Container[N]" array-of-structs line, and grow the textbox to fit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

$CR = [string][char]13

# 1) Remove the leading "1111...1" paragraph along with the blank
#    paragraph that followed it (64 '1' chars + its paragraph mark,
#    plus the following empty paragraph's mark = 66 characters).
$tr.Characters(1, 66).Delete()

# 2) Insert the new note paragraph (with a following blank paragraph)
#    right before the "Container[N]" line.
$full = $tr.Text
$idx = $full.IndexOf("Container[N]")
$insertPoint = $tr.Characters($idx + 1, 1)
$insertPoint.InsertBefore("I should first profile this in another project" + $CR + $CR)

# 3) Grow the textbox height to match the extra line of text.
$shp.Height = 4293483
